$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Part 1: the "TUE Dec 12" / " 09:42:02 PST 2017" timestamp was originally
# split across two runs; the edit merges it back into a single run. A
# Find/Replace over the whole (identically formatted) phrase collapses the
# two runs Word had split it into.
# ---------------------------------------------------------------------------
$null = $d.Content.Find.Execute(
    "TUE Dec 12 09:42:02 PST 2017", $false, $false, $false, $false, $false,
    $true, 1, $false, "TUE Dec 12 09:42:02 PST 2017", 2)

# ---------------------------------------------------------------------------
# Part 2: a new purchase record (MON Dec 18) is appended right after the
# final existing record in the document (the one ending "...Amount Received
# mode - CASH"), and before the long run of trailing blank paragraphs.
# ---------------------------------------------------------------------------

# Locate the paragraph that holds the final "- CASH" line (there are several
# "- CASH" records earlier in the document; we need the very last one).
$count = $d.Paragraphs.Count
$targetIndex = -1
for ($i = $count; $i -ge 1; $i--) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -match "- CASH") {
        $targetIndex = $i
        break
    }
}
if ($targetIndex -eq -1) {
    throw "Could not locate the '- CASH' paragraph to anchor the insertion."
}

$anchorRange = $d.Paragraphs.Item($targetIndex).Range
$insertPoint = $d.Range($anchorRange.End, $anchorRange.End)

# Word/Courier-New "PlainText" paragraph properties shared by every line of
# the record (the bold "Amount balance" line also sets <w:b/>).
$pPrPlain = '<w:pPr><w:pStyle w:val="PlainText"/><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/></w:rPr></w:pPr>'
$rPrPlain = '<w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/></w:rPr>'
$pPrBold  = '<w:pPr><w:pStyle w:val="PlainText"/><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/><w:b/></w:rPr></w:pPr>'
$rPrBold  = '<w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/><w:b/></w:rPr>'

function New-PlainRun([string]$rPr, [string]$text, [bool]$preserveSpace) {
    if ($text -eq $null) {
        return "<w:r>$rPr<w:tab/></w:r>"
    }
    if ($preserveSpace) {
        return "<w:r>$rPr<w:t xml:space=`"preserve`">$text</w:t></w:r>"
    }
    return "<w:r>$rPr<w:t>$text</w:t></w:r>"
}

function New-EmptyParagraph {
    return "<w:p>$pPrPlain</w:p>"
}

function New-DateParagraph([string]$day, [string]$rest) {
    $r1 = New-PlainRun $rPrPlain $day $false
    $r2 = New-PlainRun $rPrPlain $rest $true
    return "<w:p>$pPrPlain$r1$r2</w:p>"
}

function New-FieldParagraph([string]$label, [int]$tabCount, [string]$value, [bool]$bold) {
    $pPr = $(if ($bold) { $pPrBold } else { $pPrPlain })
    $rPr = $(if ($bold) { $rPrBold } else { $rPrPlain })
    $runs = New-PlainRun $rPr $label $false
    for ($n = 1; $n -lt $tabCount; $n++) {
        $runs += New-PlainRun $rPr $null $false
    }
    $runs += "<w:r>$rPr<w:tab/><w:t>$value</w:t></w:r>"
    return "<w:p>$pPr$runs</w:p>"
}

function New-TextParagraph([string]$text) {
    $r = New-PlainRun $rPrPlain $text $false
    return "<w:p>$pPrPlain$r</w:p>"
}

$paragraphsXml = ""
$paragraphsXml += New-EmptyParagraph
$paragraphsXml += New-DateParagraph "MON Dec 18" " 11:00:06 PST 2017"
$paragraphsXml += New-FieldParagraph "Person Name" 4 "- TK" $false
$paragraphsXml += New-FieldParagraph "Bill number" 4 "- 2031" $false
$paragraphsXml += New-TextParagraph "---------------------------------------------------------------"
$paragraphsXml += New-FieldParagraph "Item Name" 4 "- CARROT EVE" $false
$paragraphsXml += New-FieldParagraph "Number of Pockets" 3 "- 1" $false
$paragraphsXml += New-FieldParagraph "Number of KGs" 3 "- 81" $false
$paragraphsXml += New-FieldParagraph "Rate" 5 "- 40" $false
$paragraphsXml += New-FieldParagraph "Total Price" 4 "- 3240.0" $false
$paragraphsXml += New-FieldParagraph "Amount balance" 3 "- 8558.0" $true
$paragraphsXml += New-EmptyParagraph
$paragraphsXml += New-EmptyParagraph

$xml = '<?xml version="1.0" standalone="yes"?>' +
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData>' +
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:body>' + $paragraphsXml + '</w:body>' +
    '</w:document>' +
    '</pkg:xmlData></pkg:part></pkg:package>'

$insertPoint.InsertXML($xml)
